# Add final times / tasks to the time accounting sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - 2021-03-?? (row date index 24)
$ws.Range("F9").Value = "Design "
$ws.Range("N9").Value = 0.54166666666666663   # 13:00
$ws.Range("O9").Value = 0.58333333333333337   # 14:00

# Row 10 - 44266
$ws.Range("F10").Value = "Framework"
$ws.Range("H10").Value = 0.58333333333333337  # 14:00
$ws.Range("I10").Value = 0.79166666666666663  # 19:00

# Row 11 - 44267
$ws.Range("F11").Value = "Framework"
$ws.Range("H11").Value = 0.375                # 9:00
$ws.Range("I11").Value = 0.54166666666666663  # 13:00
$ws.Range("N11").Value = 0.70833333333333337  # 17:00
$ws.Range("O11").Value = 0.75                 # 18:00

# Row 12 - 44268
$ws.Range("F12").Value = "Functionality "
$ws.Range("H12").Value = 0.45833333333333331  # 11:00
$ws.Range("I12").Value = 0.70833333333333337  # 17:00
$ws.Range("N12").Value = 0.5                  # 12:00
$ws.Range("O12").Value = 0.875                # 21:00

# Row 13 - 44269
$ws.Range("B13").Value = "Groupme Conversation"
$ws.Range("D13").Value = 0.5                  # 12:00
$ws.Range("E13").Value = 0.5625                # 13:30
$ws.Range("F13").Value = "Functionality "
$ws.Range("H13").Value = 0.41666666666666669  # 10:00
$ws.Range("I13").Value = 0.625                # 15:00
$ws.Range("N13").Value = 0.54166666666666663  # 13:00
$ws.Range("O13").Value = 0.625                # 15:00
